$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 54.53585066666667
$ws.Cells.Item(2, 8).Value = 163.607552
$ws.Cells.Item(2, 9).Value = 0.3031388658437607
$ws.Cells.Item(2, 10).Value = 0.3031388658437607
$ws.Cells.Item(2, 13).Value = 3.833884
$ws.Cells.Item(2, 14).Value = 11.501652
$ws.Cells.Item(2, 15).Value = 0.06481711538755341
$ws.Cells.Item(2, 16).Value = 0.06481711538755341
$ws.Cells.Item(2, 17).Value = 209.0841252973227
$ws.Cells.Item(2, 18).Value = 1881.757127675904
$ws.Cells.Item(2, 19).Value = 0.01964858684584711
$ws.Cells.Item(2, 20).Value = 0.01964858684584711
$ws.Cells.Item(3, 7).Value = 54.53585066666667
$ws.Cells.Item(3, 8).Value = 163.607552
$ws.Cells.Item(3, 9).Value = 0.3031388658437607
$ws.Cells.Item(3, 10).Value = 0.3031388658437607
$ws.Cells.Item(3, 15).Value = 0.6503392567461023
$ws.Cells.Item(3, 16).Value = 0.6503392567461023
$ws.Cells.Item(3, 17).Value = 2097.835021355805
$ws.Cells.Item(3, 18).Value = 18880.51519220224
$ws.Cells.Item(3, 19).Value = 0.1971431047036877
$ws.Cells.Item(3, 20).Value = 0.1971431047036877
$ws.Cells.Item(4, 7).Value = 54.53585066666667
$ws.Cells.Item(4, 8).Value = 163.607552
$ws.Cells.Item(4, 9).Value = 0.3031388658437607
$ws.Cells.Item(4, 10).Value = 0.3031388658437607
$ws.Cells.Item(4, 13).Value = 9.801416999999999
$ws.Cells.Item(4, 14).Value = 29.404251
$ws.Cells.Item(4, 15).Value = 0.1657065202417516
$ws.Cells.Item(4, 16).Value = 0.1657065202417516
$ws.Cells.Item(4, 17).Value = 534.5286138337279
$ws.Cells.Item(4, 18).Value = 4810.757524503551
$ws.Cells.Item(4, 19).Value = 0.05023208660900074
$ws.Cells.Item(4, 20).Value = 0.05023208660900076
$ws.Cells.Item(5, 7).Value = 54.53585066666667
$ws.Cells.Item(5, 8).Value = 163.607552
$ws.Cells.Item(5, 9).Value = 0.3031388658437607
$ws.Cells.Item(5, 10).Value = 0.3031388658437607
$ws.Cells.Item(5, 13).Value = 7.046871
$ws.Cells.Item(5, 14).Value = 21.140613
$ws.Cells.Item(5, 15).Value = 0.1191371076245927
$ws.Cells.Item(5, 16).Value = 0.1191371076245927
$ws.Cells.Item(5, 17).Value = 384.3071045232641
$ws.Cells.Item(5, 18).Value = 3458.763940709376
$ws.Cells.Item(5, 19).Value = 0.03611508768522508
$ws.Cells.Item(5, 20).Value = 0.03611508768522509
$ws.Cells.Item(6, 9).Value = 0.1026363515063155
$ws.Cells.Item(6, 10).Value = 0.1026363515063155
$ws.Cells.Item(6, 13).Value = 3.833884
$ws.Cells.Item(6, 14).Value = 11.501652
$ws.Cells.Item(6, 15).Value = 0.06481711538755341
$ws.Cells.Item(6, 16).Value = 0.06481711538755341
$ws.Cells.Item(6, 17).Value = 70.79142332566133
$ws.Cells.Item(6, 18).Value = 637.1228099309519
$ws.Cells.Item(6, 19).Value = 0.006652592238542343
$ws.Cells.Item(6, 20).Value = 0.006652592238542344
$ws.Cells.Item(7, 9).Value = 0.1026363515063155
$ws.Cells.Item(7, 10).Value = 0.1026363515063155
$ws.Cells.Item(7, 15).Value = 0.6503392567461023
$ws.Cells.Item(7, 16).Value = 0.6503392567461023
$ws.Cells.Item(7, 19).Value = 0.06674844855374891
$ws.Cells.Item(7, 20).Value = 0.06674844855374892
$ws.Cells.Item(8, 9).Value = 0.1026363515063155
$ws.Cells.Item(8, 10).Value = 0.1026363515063155
$ws.Cells.Item(8, 13).Value = 9.801416999999999
$ws.Cells.Item(8, 14).Value = 29.404251
$ws.Cells.Item(8, 15).Value = 0.1657065202417516
$ws.Cells.Item(8, 16).Value = 0.1657065202417516
$ws.Cells.Item(8, 17).Value = 180.979982711614
$ws.Cells.Item(8, 18).Value = 1628.819844404526
$ws.Cells.Item(8, 19).Value = 0.0170075126584208
$ws.Cells.Item(8, 20).Value = 0.0170075126584208
$ws.Cells.Item(9, 9).Value = 0.1026363515063155
$ws.Cells.Item(9, 10).Value = 0.1026363515063155
$ws.Cells.Item(9, 13).Value = 7.046871
$ws.Cells.Item(9, 14).Value = 21.140613
$ws.Cells.Item(9, 15).Value = 0.1191371076245927
$ws.Cells.Item(9, 16).Value = 0.1191371076245927
$ws.Cells.Item(9, 17).Value = 130.118185130882
$ws.Cells.Item(9, 18).Value = 1171.063666177938
$ws.Cells.Item(9, 19).Value = 0.01222779805560343
$ws.Cells.Item(9, 20).Value = 0.01222779805560344
$ws.Cells.Item(10, 7).Value = 12.55635966666667
$ws.Cells.Item(10, 8).Value = 37.669079
$ws.Cells.Item(10, 9).Value = 0.06979483370938171
$ws.Cells.Item(10, 10).Value = 0.06979483370938172
$ws.Cells.Item(10, 13).Value = 3.833884
$ws.Cells.Item(10, 14).Value = 11.501652
$ws.Cells.Item(10, 15).Value = 0.06481711538755341
$ws.Cells.Item(10, 16).Value = 0.06481711538755341
$ws.Cells.Item(10, 17).Value = 48.13962642427866
$ws.Cells.Item(10, 18).Value = 433.256637818508
$ws.Cells.Item(10, 19).Value = 0.004523899789996097
$ws.Cells.Item(10, 20).Value = 0.004523899789996098
$ws.Cells.Item(11, 7).Value = 12.55635966666667
$ws.Cells.Item(11, 8).Value = 37.669079
$ws.Cells.Item(11, 9).Value = 0.06979483370938171
$ws.Cells.Item(11, 10).Value = 0.06979483370938172
$ws.Cells.Item(11, 15).Value = 0.6503392567461023
$ws.Cells.Item(11, 16).Value = 0.6503392567461023
$ws.Cells.Item(11, 17).Value = 483.0065127337061
$ws.Cells.Item(11, 18).Value = 4347.058614603354
$ws.Cells.Item(11, 19).Value = 0.04539032027927711
$ws.Cells.Item(11, 20).Value = 0.04539032027927711
$ws.Cells.Item(12, 7).Value = 12.55635966666667
$ws.Cells.Item(12, 8).Value = 37.669079
$ws.Cells.Item(12, 9).Value = 0.06979483370938171
$ws.Cells.Item(12, 10).Value = 0.06979483370938172
$ws.Cells.Item(12, 13).Value = 9.801416999999999
$ws.Cells.Item(12, 14).Value = 29.404251
$ws.Cells.Item(12, 15).Value = 0.1657065202417516
$ws.Cells.Item(12, 16).Value = 0.1657065202417516
$ws.Cells.Item(12, 17).Value = 123.070117094981
$ws.Cells.Item(12, 18).Value = 1107.631053854829
$ws.Cells.Item(12, 19).Value = 0.01156545902483335
$ws.Cells.Item(12, 20).Value = 0.01156545902483335
$ws.Cells.Item(13, 7).Value = 12.55635966666667
$ws.Cells.Item(13, 8).Value = 37.669079
$ws.Cells.Item(13, 9).Value = 0.06979483370938171
$ws.Cells.Item(13, 10).Value = 0.06979483370938172
$ws.Cells.Item(13, 13).Value = 7.046871
$ws.Cells.Item(13, 14).Value = 21.140613
$ws.Cells.Item(13, 15).Value = 0.1191371076245927
$ws.Cells.Item(13, 16).Value = 0.1191371076245927
$ws.Cells.Item(13, 17).Value = 88.48304680060299
$ws.Cells.Item(13, 18).Value = 796.3474212054269
$ws.Cells.Item(13, 19).Value = 0.008315154615275159
$ws.Cells.Item(13, 20).Value = 0.008315154615275161
$ws.Cells.Item(14, 7).Value = 94.34696966666667
$ws.Cells.Item(14, 8).Value = 283.040909
$ws.Cells.Item(14, 9).Value = 0.524429948940542
$ws.Cells.Item(14, 10).Value = 0.5244299489405421
$ws.Cells.Item(14, 13).Value = 3.833884
$ws.Cells.Item(14, 14).Value = 11.501652
$ws.Cells.Item(14, 15).Value = 0.06481711538755341
$ws.Cells.Item(14, 16).Value = 0.06481711538755341
$ws.Cells.Item(14, 17).Value = 361.7153374535187
$ws.Cells.Item(14, 18).Value = 3255.438037081668
$ws.Cells.Item(14, 19).Value = 0.03399203651316785
$ws.Cells.Item(14, 20).Value = 0.03399203651316786
$ws.Cells.Item(15, 7).Value = 94.34696966666667
$ws.Cells.Item(15, 8).Value = 283.040909
$ws.Cells.Item(15, 9).Value = 0.524429948940542
$ws.Cells.Item(15, 10).Value = 0.5244299489405421
$ws.Cells.Item(15, 15).Value = 0.6503392567461023
$ws.Cells.Item(15, 16).Value = 0.6503392567461023
$ws.Cells.Item(15, 17).Value = 3629.252587170189
$ws.Cells.Item(15, 18).Value = 32663.2732845317
$ws.Cells.Item(15, 19).Value = 0.3410573832093884
$ws.Cells.Item(15, 20).Value = 0.3410573832093885
$ws.Cells.Item(16, 7).Value = 94.34696966666667
$ws.Cells.Item(16, 8).Value = 283.040909
$ws.Cells.Item(16, 9).Value = 0.524429948940542
$ws.Cells.Item(16, 10).Value = 0.5244299489405421
$ws.Cells.Item(16, 13).Value = 9.801416999999999
$ws.Cells.Item(16, 14).Value = 29.404251
$ws.Cells.Item(16, 15).Value = 0.1657065202417516
$ws.Cells.Item(16, 16).Value = 0.1657065202417516
$ws.Cells.Item(16, 17).Value = 924.7339923893509
$ws.Cells.Item(16, 18).Value = 8322.605931504158
$ws.Cells.Item(16, 19).Value = 0.08690146194949666
$ws.Cells.Item(16, 20).Value = 0.0869014619494967
$ws.Cells.Item(17, 7).Value = 94.34696966666667
$ws.Cells.Item(17, 8).Value = 283.040909
$ws.Cells.Item(17, 9).Value = 0.524429948940542
$ws.Cells.Item(17, 10).Value = 0.5244299489405421
$ws.Cells.Item(17, 13).Value = 7.046871
$ws.Cells.Item(17, 14).Value = 21.140613
$ws.Cells.Item(17, 15).Value = 0.1191371076245927
$ws.Cells.Item(17, 16).Value = 0.1191371076245927
$ws.Cells.Item(17, 17).Value = 664.850924481913
$ws.Cells.Item(17, 18).Value = 5983.658320337217
$ws.Cells.Item(17, 19).Value = 0.06247906726848899
$ws.Cells.Item(17, 20).Value = 0.06247906726848901
